$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D18/E18 values (plain numbers, no formulas)
$ws.Range("D18").Value = 30000
$ws.Range("E18").Value = 20000

# Update D19/E19 with formulas referencing row above, doubled
$ws.Range("D19").Formula = "=D18*2"
$ws.Range("E19").Formula = "=E18*2"

# Update D20/E20 with formulas referencing row above, doubled
$ws.Range("D20").Formula = "=D19*2"
$ws.Range("E20").Formula = "=E19*2"

# Update the sheet view state: topLeftCell and selected cell
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D16").Select()
